$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) and, for row 34, column E (dS0) to repulled / recalculated values.
$ws.Range("F9").Value = -1
$ws.Range("F14").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 4
$ws.Range("F21").Value = 2
$ws.Range("F27").Value = -1
$ws.Range("F30").Value = 3
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 4
$ws.Range("F39").Value = 2
$ws.Range("F41").Value = -2
$ws.Range("F43").Value = 0
$ws.Range("F45").Value = 2
$ws.Range("F47").Value = -2
$ws.Range("F48").Value = 2
$ws.Range("F49").Value = 4
$ws.Range("F53").Value = 1
$ws.Range("F55").Value = -2
$ws.Range("F60").Value = -1
$ws.Range("F62").Value = -3
$ws.Range("F69").Value = 2
$ws.Range("F71").Value = 4
